$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new Rank (column A) values for rows 47-51 that were previously missing
$ws.Range("A47").Value = 46
$ws.Range("A48").Value = 47
$ws.Range("A49").Value = 48
$ws.Range("A50").Value = 49
$ws.Range("A51").Value = 50

# Update recalculated market data values (columns J and L)
$ws.Range("J2").Value = 25.64
$ws.Range("J3").Value = 1321
$ws.Range("J4").Value = 6.54
$ws.Range("J5").Value = 33.67
$ws.Range("J6").Value = 60.86
$ws.Range("J7").Value = 29.4
$ws.Range("J8").Value = 455.9
$ws.Range("J9").Value = 732.9
$ws.Range("L9").Value = 1.1852
$ws.Range("J10").Value = 56.9
$ws.Range("J11").Value = 116.6
$ws.Range("J12").Value = 15.04
$ws.Range("J13").Value = 78.22
$ws.Range("J14").Value = 83.5
$ws.Range("J15").Value = 65.52
$ws.Range("J16").Value = 91.62
$ws.Range("J17").Value = 21.505
$ws.Range("J18").Value = 10.83
$ws.Range("J19").Value = 31.08
$ws.Range("J21").Value = 4.173
$ws.Range("J22").Value = 57.54
$ws.Range("J23").Value = 23.7
$ws.Range("J24").Value = 54.63
$ws.Range("J25").Value = 35.65
$ws.Range("L25").Value = 1.1852
$ws.Range("J26").Value = 39.44
$ws.Range("L26").Value = 1.1852
$ws.Range("J27").Value = 30.56
$ws.Range("J28").Value = 24.6
$ws.Range("J29").Value = 52.52
$ws.Range("J30").Value = 21.04
$ws.Range("J31").Value = 110.35
$ws.Range("J32").Value = 49.04
$ws.Range("J33").Value = 229.6
$ws.Range("L33").Value = 1.1852
$ws.Range("J34").Value = 62.7
$ws.Range("L34").Value = 1.1852
$ws.Range("J35").Value = 57.85
$ws.Range("J36").Value = 630.6
$ws.Range("J37").Value = 12.18
$ws.Range("J38").Value = 74.40000000000001
$ws.Range("J39").Value = 162.2
$ws.Range("J40").Value = 95.34999999999999
$ws.Range("J41").Value = 74.3
$ws.Range("L41").Value = 1.1852
$ws.Range("J42").Value = 100.85
$ws.Range("J43").Value = 8100
$ws.Range("J44").Value = 64.45
$ws.Range("J45").Value = 66.40000000000001
$ws.Range("L45").Value = 1.1852
$ws.Range("J46").Value = 249.6
$ws.Range("J48").Value = 71.72
$ws.Range("J49").Value = 203.9
$ws.Range("J50").Value = 12.78
$ws.Range("J51").Value = 21.04
